# YPS Startup Expenses.xlsx - update header of first section
# (commit: "Updated header of first section in YPS Startup Expenses.xlsx")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Startup Expenses")

# A7 previously read "Owners' Investment (name & % ownership)" (wrapped over two
# lines, hence the taller row). Shorten it to "Owners' Investment (names)".
$ws.Range("A7").Value = "Owners' Investment (names)"

# The shorter text no longer needs to wrap across two lines, so the explicit
# row height goes away again (back to the sheet's default row height).
$ws.Rows(7).AutoFit() | Out-Null

# Move the saved cursor/selection from where editing finished (E76, scrolled
# down to row 88) back up to the top of the sheet at A8.
$ws.Range("A8").Select() | Out-Null
